# #5: property aircraft done
# Fix the property_category column values on the "building" (建物) and
# "car" (汽車) sheets, which had been mistakenly left as "land".

$wb = $excel.ActiveWorkbook

# 建物 sheet: row 2, column I (property_category) -> "building"
$wsBuilding = $wb.Worksheets.Item("建物")
$wsBuilding.Range("I2").Value = "building"

# 汽車 sheet: rows 2-4, column H (property_category) -> "car"
$wsCar = $wb.Worksheets.Item("汽車")
$wsCar.Range("H2").Value = "car"
$wsCar.Range("H3").Value = "car"
$wsCar.Range("H4").Value = "car"
